# "Add contests part admin"
#
# Adds a new "10-AdminLogin" worksheet (admin login credentials, mirroring
# the existing login-data sheets) at the end of the workbook, and makes it
# the active/selected sheet - which in turn de-selects the sheet that used
# to be active ("5-LoginUserTest").

$wb = $excel.ActiveWorkbook

# --- 1. Update the previously-active sheet's cursor position -------------
# "5-LoginUserTest" was the active tab before this edit; once the new sheet
# becomes active it keeps a plain (non-selected) sheetView with its cursor
# left on A2.
$loginSheet = $wb.Worksheets.Item("5-LoginUserTest")
$loginSheet.Activate()
$loginSheet.Range("A2").Select() | Out-Null

# --- 2. Add the new worksheet at the end of the workbook -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "10-AdminLogin"

# --- 3. Fill in the admin login data --------------------------------------
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 1

$ws.Range("A2").Value = "emailAdmin"
$ws.Range("B2").Value = "ahmadali@hawyah.com"

$ws.Range("A3").Value = "password"
$ws.Range("B3").Value = "11111111a"
$ws.Range("B3").NumberFormat = "0;[Red]0"

# --- 4. Hyperlink the admin email address ---------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:ahmadali@hawyah.com") | Out-Null
$ws.Range("B2").Style = "Normal"

# --- 5. Column widths (matching the other login sheets' autofit look) ----
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(2).ColumnWidth = 20

# --- 6. Page setup ----------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- 7. Leave the cursor where the author left it, and make the sheet active
$ws.Range("B4").Select() | Out-Null
$ws.Activate()
